# Singapore GP comparison table fixes: "missing times" correction.
# 1) Update cumulative_time_error / gap_error values for rows 2-19
#    (Lando Norris .. Daniel Ricciardo).
# 2) Swap the order of the final two rows (Alexander Albon / Kevin
#    Magnussen) so Kevin Magnussen now precedes Alexander Albon.

$d = $word.ActiveDocument

function Replace-Value($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# cumulative_time_error / gap_error corrections, row by row.
Replace-Value "197.492" "199.587"

Replace-Value "203.272" "205.933"
Replace-Value "5.780"   "6.345"

Replace-Value "199.402" "203.962"
Replace-Value "1.909"   "4.374"

Replace-Value "196.924" "202.350"
Replace-Value "-0.568"  "2.763"

Replace-Value "201.099" "205.045"
Replace-Value "3.606"   "5.458"

Replace-Value "199.128" "202.103"
Replace-Value "1.636"   "2.515"

Replace-Value "190.335" "196.784"
Replace-Value "-7.157"  "-2.803"

Replace-Value "290.508" "295.379"
Replace-Value "-4.906"  "-2.130"

Replace-Value "285.656" "291.239"
Replace-Value "-9.758"  "-6.270"

Replace-Value "295.408" "299.533"
Replace-Value "-0.006"  "2.024"

Replace-Value "296.117" "302.610"
Replace-Value "0.702"   "5.101"

Replace-Value "294.430" "300.666"
Replace-Value "-0.985"  "3.157"

Replace-Value "281.072" "289.046"
Replace-Value "-14.342" "-8.463"

Replace-Value "286.718" "293.556"
Replace-Value "-8.697"  "-3.953"

Replace-Value "278.677" "286.325"
Replace-Value "-16.737" "-11.185"

Replace-Value "286.119" "294.454"
Replace-Value "-9.296"  "-3.055"

Replace-Value "289.183" "297.314"
Replace-Value "-6.232"  "-0.196"

Replace-Value "317.749" "324.651"
Replace-Value "22.334"  "27.142"

# Swap the last two data rows (Alexander Albon <-> Kevin Magnussen) by
# exchanging their cell contents, which reproduces the row reordering
# seen in the diff without relying on row-move/insert/delete APIs.
$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count
$albonRow = $rowCount - 1
$magnussenRow = $rowCount
$colCount = $t.Columns.Count

$albonVals = @()
for ($col = 1; $col -le $colCount; $col++) {
    $albonVals += $t.Cell($albonRow, $col).Range.Text
}
$magnussenVals = @()
for ($col = 1; $col -le $colCount; $col++) {
    $magnussenVals += $t.Cell($magnussenRow, $col).Range.Text
}

for ($col = 1; $col -le $colCount; $col++) {
    $t.Cell($albonRow, $col).Range.Text = $magnussenVals[$col - 1]
}
for ($col = 1; $col -le $colCount; $col++) {
    $t.Cell($magnussenRow, $col).Range.Text = $albonVals[$col - 1]
}
